$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "plain" numeric-looking string (e.g. "375.91")
# must be forced to text, otherwise Excel auto-converts them to numbers.
# Force via NumberFormat "@" then restore formatting so no stray style sticks.
$textForceCells = @(
    "D5",
    "D6",
    "D9",
    "D10",
    "D12",
    "D14",
    "D15",
    "D17",
    "D19",
    "D20",
    "D21",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D38",
    "D39",
    "D40",
    "D43",
    "D44",
    "D46",
    "D48",
    "D50",
    "D51",
)
foreach ($ref in $textForceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply all literal cell values
$ws.Range("D2").Value = "51.128.25"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "2.939.10"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "375.91"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").Value = "102.48"
$ws.Range("E6").Value = "  -3.52%  "
$ws.Range("E7").Value = "  -1.78%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.583"
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("D10").Value = "36.72"
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").Value = "0.0836"
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").Value = "3.407.20"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").Value = "17.93"
$ws.Range("E14").Value = "  -4.10%  "
$ws.Range("D15").Value = "7.34"
$ws.Range("E15").Value = "  -2.34%  "
$ws.Range("D16").Value = "2.934.86"
$ws.Range("D17").Value = "0.975"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "51.087.11"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").Value = "3.15"
$ws.Range("E19").Value = "  -7.03%  "
$ws.Range("D20").Value = "7.11"
$ws.Range("E20").Value = "  -4.15%  "
$ws.Range("D21").Value = "12.57"
$ws.Range("E21").Value = "  -3.50%  "
$ws.Range("D22").Value = "0.0₃0953"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").Value = "263.05"
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").Value = "68.19"
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("D25").Value = "2.87"
$ws.Range("E25").Value = "  +2.52%  "
$ws.Range("B26").Value = "Filecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D26").Value = "8.12"
$ws.Range("E26").Value = "  +8.87%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "7.84"
$ws.Range("E27").Value = "  +8.40%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("E29").Value = "  -2.67%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "0.113"
$ws.Range("E30").Value = "  +4.13%  "
$ws.Range("D31").Value = "25.66"
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("D32").Value = "9.86"
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("B33").Value = "VeChain"
$ws.Range("C33").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D33").Value = "0.0457"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("B34").Value = "OKB"
$ws.Range("C34").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D34").Value = "50.92"
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "34.14"
$ws.Range("E35").Value = "  -2.10%  "
$ws.Range("E36").Value = "  -3.74%  "
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "2.97"
$ws.Range("E38").Value = "  -4.17%  "
$ws.Range("D39").Value = "2.57"
$ws.Range("E39").Value = "  -2.26%  "
$ws.Range("D40").Value = "16.44"
$ws.Range("E40").Value = "  -5.77%  "
$ws.Range("E41").Value = "  -1.34%  "
$ws.Range("E42").Value = "  -3.82%  "
$ws.Range("D43").Value = "121.83"
$ws.Range("E43").Value = "  -1.56%  "
$ws.Range("D44").Value = "21.07"
$ws.Range("E44").Value = "  -5.05%  "
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("D46").Value = "0.272"
$ws.Range("E46").Value = "  -3.13%  "
$ws.Range("E47").Value = "  -2.53%  "
$ws.Range("D48").Value = "3.22"
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("D49").Value = "1.997.16"
$ws.Range("E49").Value = "  -2.51%  "
$ws.Range("D50").Value = "0.0352"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").Value = "5.03"
$ws.Range("E51").Value = "  -3.00%  "

# Remove the temporary text-format styling so cells keep their original (default) style
foreach ($ref in $textForceCells) {
    $ws.Range($ref).ClearFormats()
}